# Fruta / hortaliza, semanal
#
# This weekly refresh prepends a new week's worth of data (two rows: a
# "Primera" and "Segunda" quality record, both dated 44483) to the
# existing daily log, pushing all of the previously-recorded rows (old
# rows 68..157) down by two rows (they land at rows 70..159) without any
# other change to their contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows above the current row 68; Excel shifts
# everything from row 68 down onward by two rows automatically, which
# also grows the used range from A1:R157 to A1:R159.
$ws.Rows("68:69").Insert()

# Populate the first of the two newly inserted rows.
$ws.Range("A68").Value = 9
$ws.Range("B68").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C68").Value = "Metropolitana"
$ws.Range("D68").Value = 44483
$ws.Range("E68").Value = 13
$ws.Range("F68").Value = 100112043
$ws.Range("G68").Value = "Pepino ensalada"
$ws.Range("H68").Value = "Sin especificar"
$ws.Range("I68").Value = "Primera"
$ws.Range("J68").Value = 79
$ws.Range("K68").Value = 12000
$ws.Range("L68").Value = 13000
$ws.Range("M68").Value = 12494
$ws.Range("N68").Value = "`$/caja 60 unidades"
$ws.Range("O68").Value = "Región de Arica y Parinacota"
$ws.Range("P68").Value = 208
$ws.Range("Q68").Value = 60
$ws.Range("R68").Value = "Hortaliza"

# Populate the second of the two newly inserted rows.
$ws.Range("A69").Value = 9
$ws.Range("B69").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C69").Value = "Metropolitana"
$ws.Range("D69").Value = 44483
$ws.Range("E69").Value = 13
$ws.Range("F69").Value = 100112043
$ws.Range("G69").Value = "Pepino ensalada"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Segunda"
$ws.Range("J69").Value = 43
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = 10000
$ws.Range("N69").Value = "`$/caja 100 unidades"
$ws.Range("O69").Value = "Región de Arica y Parinacota"
$ws.Range("P69").Value = 100
$ws.Range("Q69").Value = 100
$ws.Range("R69").Value = "Hortaliza"
